# "Generate Report for Handoff"
# The localization status report was regenerated: the in-flight items moved
# from "In Translation" to "Ready for handoff", the handoff timestamps were
# refreshed, and the Status columns were widened to fit the new (longer)
# status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value     = "Ready for handoff"   # Status
$wsDeDe.Range("C2").Value     = "Ready for handoff"   # Status

# --- Refreshed handoff / report-generation timestamps ----------------------
$wsOverview.Range("G2").Value = "2016-08-16 16:35:30" # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value     = "2016-08-16 16:35:30" # Latest Handoff Datetime (de-de)
$wsZhCn.Range("H2").Value     = "2016-08-16 16:35:25" # Latest Handoff Datetime (zh-cn)

# --- Widen the Status columns to fit "Ready for handoff" --------------------
$newStatusWidth = 16.333333333333332
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth  # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth  # F: de-de status
$wsZhCn.Columns.Item(3).ColumnWidth     = $newStatusWidth  # C: Status
$wsDeDe.Columns.Item(3).ColumnWidth     = $newStatusWidth  # C: Status
